# LV_TMTI0056883_VerifyNewJobTypeIsForFVAOnly.xlsx
# "Merge - Opp Test Data, ENg Detail, Add Counterparty - 10 Oct 2025"
#
# - Users sheet, cell B2: "Liz Hedgcock" -> "Blaise Brunda"
# - Users sheet: active selection moves from B2 to E19

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Update the user name in B2 (shared string table update)
$ws.Range("B2").Value = "Blaise Brunda"

# Update the saved selection/active cell for the sheet
$ws.Activate()
$ws.Range("E19").Select() | Out-Null
